$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.756.05'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.353.55'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.28'
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.666'
$ws.Range("E6").Value = '  -1.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.65'
$ws.Range("E7").Value = '  +1.50%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.595'
$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("E10").Value = '  +2.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.73'
$ws.Range("E11").Value = '  +6.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '37.12'
$ws.Range("E12").Value = '  +12.55%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.27'
$ws.Range("E14").Value = '  -0.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.702.25'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.38'
$ws.Range("E16").Value = '  -0.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.938'
$ws.Range("E17").Value = '  +4.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.359.48'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.711.52'
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("E20").Value = '  +2.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.59'
$ws.Range("E21").Value = '  -3.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.48'
$ws.Range("E22").Value = '  +0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.80'
$ws.Range("E23").Value = '  -1.19%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.76'
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -6.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.51'
$ws.Range("E27").Value = '  +1.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.70'
$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("E29").Value = '  +1.41%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.81'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.23'
$ws.Range("E31").Value = '  -1.47%  '

$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.135'
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0764'
$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.57'
$ws.Range("E35").Value = '  +0.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.19'
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.80'
$ws.Range("E37").Value = '  +1.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.63'
$ws.Range("E38").Value = '  +7.00%  '

$ws.Range("E39").Value = '  +1.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0282'
$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.22'
$ws.Range("E41").Value = '  +12.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.45'
$ws.Range("E42").Value = '  +17.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.23'
$ws.Range("E43").Value = '  +7.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.109'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.19'
$ws.Range("E45").Value = '  +2.05%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.202'
$ws.Range("E46").Value = '  -1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.55'
$ws.Range("E47").Value = '  +3.67%  '

$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.25'
$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("B49").Value = 'BinanceUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.15'
$ws.Range("E51").Value = '  -1.51%  '

Write-Host "Update complete"